# Quarterly financials update: insert two new quarter columns (D, E) before
# the existing data, shifting prior quarters right, and refresh the values
# for the sheet (including a couple of data corrections to older quarters).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new blank columns at D (existing D:K data shifts to F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# 2. Copy number/date formatting from column F (the shifted former column D)
#    into the two newly inserted, currently unformatted columns D:E so the
#    new quarter columns render the same as the others (date row vs. data rows).
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$cellValues = @{
    "D7" = 43465
    "E7" = 43373
    "D8" = 464500
    "E8" = 238300
    "J8" = 830100
    "D9" = 218800
    "E9" = 51300
    "J9" = 345700
    "D10" = 245700
    "E10" = 187000
    "J10" = 484400
    "D12" = "NA"
    "E12" = "NA"
    "D13" = 0
    "E13" = 0
    "D14" = 0
    "E14" = 0
    "D15" = 41100
    "E15" = 40300
    "J15" = 74900
    "D17" = 383600
    "E17" = 202000
    "H17" = 368600
    "I17" = 202000
    "J17" = 652300
    "D18" = 80900
    "E18" = 36300
    "H18" = 93800
    "I18" = 45100
    "J18" = 177800
    "D20" = -5100
    "E20" = -1900
    "H20" = -3500
    "I20" = -3700
    "J20" = -7300
    "D21" = 116900
    "E21" = 74700
    "J21" = 245400
    "D22" = 14600
    "E22" = 12400
    "J22" = 22800
    "D23" = 61200
    "E23" = 21900
    "J23" = 147700
    "D24" = 16500
    "E24" = 5700
    "J24" = 50600
    "D25" = 0
    "E25" = 0
    "D26" = 44700
    "E26" = 16300
    "J26" = 97100
    "D27" = 44700
    "E27" = 16300
    "J27" = 97100
    "D28" = 0
    "E28" = 0
    "D29" = "NA"
    "E29" = "NA"
    "D30" = 0
    "E30" = 0
    "D31" = 0
    "E31" = 0
    "D32" = 5100
    "E32" = 1900
    "H32" = 3500
    "I32" = 3700
    "J32" = 7300
    "D33" = 44700
    "E33" = 16300
    "J33" = 97100
    "D34" = 0
    "E34" = 0
    "D35" = 44700
    "E35" = 16300
    "J35" = 97100
    "D38" = 43465
    "E38" = 43373
    "D41" = 21300
    "E41" = 12400
    "D42" = 0
    "E42" = 0
    "D43" = 295400
    "E43" = 132400
    "D44" = 151600
    "E44" = 166800
    "D45" = 74900
    "E45" = 62800
    "D46" = 543300
    "E46" = 374500
    "D47" = 0
    "E47" = 0
    "D48" = 4283700
    "E48" = 4195900
    "D49" = 158000
    "E49" = 158000
    "D50" = 0
    "E50" = 0
    "D51" = 0
    "E51" = 0
    "D52" = 483700
    "E52" = 424600
    "D53" = 0
    "E53" = 0
    "D54" = 5468600
    "E54" = 5153000
    "D57" = 174500
    "E57" = 68300
    "D58" = 299500
    "E58" = 576000
    "D59" = 224900
    "E59" = 191700
    "D60" = 698900
    "E60" = 836000
    "D61" = 1285500
    "E61" = 893900
    "D62" = 1441600
    "E62" = 1406400
    "D63" = 0
    "E63" = 0
    "D64" = 0
    "E64" = 0
    "D65" = 0
    "E65" = 0
    "D66" = 3426000
    "E66" = 3136300
    "D68" = 0
    "E68" = 0
    "D69" = 0
    "E69" = 0
    "D70" = 0
    "E70" = 0
    "D71" = 0
    "E71" = 0
    "D72" = 320900
    "E72" = 300500
    "D73" = 0
    "E73" = 0
    "D74" = 0
    "E74" = 0
    "D75" = 0
    "E75" = 0
    "D76" = 2042700
    "E76" = 2016600
    "D77" = 0
    "E77" = 0
    "D80" = 43465
    "E80" = 43373
    "D81" = 44700
    "E81" = 16300
    "J81" = 97100
    "D83" = 41100
    "E83" = 40300
    "J83" = 74900
    "D84" = 0
    "E84" = 0
    "D85" = 0
    "E85" = 0
    "D86" = 0
    "E86" = 0
    "D87" = 0
    "E87" = 0
    "D88" = 0
    "E88" = 0
    "D89" = 30900
    "E89" = 36800
    "J89" = 279300
    "D91" = -115100
    "E91" = -103500
    "J91" = -154700
    "D92" = 0
    "E92" = 0
    "D93" = 0
    "E93" = 0
    "D94" = -115100
    "E94" = -103500
    "J94" = -154200
    "D96" = -24200
    "E96" = -24200
    "J96" = -44000
    "D97" = 0
    "E97" = 0
    "D98" = 0
    "E98" = 0
    "D99" = 0
    "E99" = 0
    "D100" = 93100
    "E100" = 66600
    "J100" = -134600
    "D101" = 0
    "E101" = 0
    "D102" = 8900
    "E102" = -200
    "J102" = -9600
}
# 3. Write the new/ corrected cell values (new quarter data for columns D
#    and E, plus a handful of restated values for older quarters).
foreach ($addr in $cellValues.Keys) {
    $ws.Range($addr).Value2 = $cellValues[$addr]
}
